$p = $ppt.ActivePresentation

# 1) Fix the body text on the "Strike similar contracts..." slide (slide 6):
#    drop the word "similar". Locate the shape by its current text so the
#    script isn't brittle to shape-ordering assumptions.
$oldText = "Strike similar contracts with Iowa to require SAT testing for admissions"
$newText = "Strike contracts with Iowa to require SAT testing for admissions"

$target = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq $oldText) {
                $target = $shape
            }
        }
    }
}
if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = $newText
}

# 2) Theme colors: restore the deck's theme color scheme to the original
#    "Default" Office colors (it currently holds the "Swiss" palette). The
#    12 theme color slots are addressed in the standard order: dk1, lt1,
#    dk2, lt2, accent1-6, hlink, folHlink. RGB() packs R,G,B (0-255) into
#    the COM BGR-ordered integer PowerPoint expects for .RGB.
function RGB([int]$r, [int]$g, [int]$b) { return $r + ($g * 256) + ($b * 65536) }

$cs = $p.SlideMaster.Theme.ThemeColorScheme
$cs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$cs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$cs.Item(3).RGB  = RGB 0x15 0x81 0x58   # dk2
$cs.Item(4).RGB  = RGB 0xF3 0xF3 0xF3   # lt2
$cs.Item(5).RGB  = RGB 0x05 0x8D 0xC7   # accent1
$cs.Item(6).RGB  = RGB 0x50 0xB4 0x32   # accent2
$cs.Item(7).RGB  = RGB 0xED 0x56 0x1B   # accent3
$cs.Item(8).RGB  = RGB 0xED 0xEF 0x00   # accent4
$cs.Item(9).RGB  = RGB 0x24 0xCB 0xE5   # accent5
$cs.Item(10).RGB = RGB 0x64 0xE5 0x72   # accent6
$cs.Item(11).RGB = RGB 0x22 0x00 0xCC   # hlink
$cs.Item(12).RGB = RGB 0x55 0x1A 0x8B   # folHlink
